$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.755.35'
$ws.Range("E2").Value = '  +7.18%  '

$ws.Range("D3").Value = '1.741.06'
$ws.Range("E3").Value = '  +3.67%  '

$cell = $ws.Range("D4")
$cell.Value = "'1.001"
$cell.Style = "Normal"
$ws.Range("E4").Value = '  -0.01%  '

$cell = $ws.Range("D5")
$cell.Value = "'335.29"
$cell.Style = "Normal"
$ws.Range("E5").Value = '  +1.74%  '

$cell = $ws.Range("D6")
$cell.Value = "'0.9985"
$cell.Style = "Normal"
$ws.Range("E6").Value = '  -0.12%  '

$cell = $ws.Range("D7")
$cell.Value = "'0.3745"
$cell.Style = "Normal"
$ws.Range("E7").Value = '  +2.30%  '

$ws.Range("B8").Value = 'OKB'
$ws.Range("C8").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$cell = $ws.Range("D8")
$cell.Value = "'48.46"
$cell.Style = "Normal"
$ws.Range("E8").Value = '  +2.80%  '

$ws.Range("B9").Value = 'Cardano'
$ws.Range("C9").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$cell = $ws.Range("D9")
$cell.Value = "'0.3396"
$cell.Style = "Normal"
$ws.Range("E9").Value = '  +4.21%  '

$cell = $ws.Range("D10")
$cell.Value = "'1.189"
$cell.Style = "Normal"
$ws.Range("E10").Value = '  +3.51%  '

$cell = $ws.Range("D11")
$cell.Value = "'0.07485"
$cell.Style = "Normal"
$ws.Range("E11").Value = '  +5.36%  '

$cell = $ws.Range("D12")
$cell.Value = "'0.9993"
$cell.Style = "Normal"
$ws.Range("E12").Value = '  +0.02%  '

$cell = $ws.Range("D13")
$cell.Value = "'6.418"
$cell.Style = "Normal"
$ws.Range("E13").Value = '  +5.06%  '

$cell = $ws.Range("D14")
$cell.Value = "'20.47"
$cell.Style = "Normal"
$ws.Range("E14").Value = '  +3.68%  '

$cell = $ws.Range("D15")
$cell.Value = "'7.073"
$cell.Style = "Normal"
$ws.Range("E15").Value = '  +6.44%  '

$ws.Range("D16").Value = '1.741.82'
$ws.Range("E16").Value = '  +3.66%  '

$cell = $ws.Range("D17")
$cell.Value = "'0.00001079"
$cell.Style = "Normal"
$ws.Range("E17").Value = '  +2.45%  '

$cell = $ws.Range("D18")
$cell.Value = "'0.06729"
$cell.Style = "Normal"
$ws.Range("E18").Value = '  +1.78%  '

$cell = $ws.Range("D19")
$cell.Value = "'82.72"
$cell.Style = "Normal"
$ws.Range("E19").Value = '  +4.52%  '

$cell = $ws.Range("D20")
$cell.Value = "'0.9979"
$cell.Style = "Normal"
$ws.Range("E20").Value = '  -0.13%  '

$cell = $ws.Range("D21")
$cell.Value = "'16.75"
$cell.Style = "Normal"
$ws.Range("E21").Value = '  +4.70%  '

$cell = $ws.Range("D22")
$cell.Value = "'6.223"
$cell.Style = "Normal"
$ws.Range("E22").Value = '  +4.79%  '

$cell = $ws.Range("D23")
$cell.Value = "'12.78"
$cell.Style = "Normal"
$ws.Range("E23").Value = '  -1.04%  '

$ws.Range("D24").Value = '26.755.58'
$ws.Range("E24").Value = '  +7.29%  '

$cell = $ws.Range("D25")
$cell.Value = "'2.452"
$cell.Style = "Normal"
$ws.Range("E25").Value = '  +0.31%  '

$cell = $ws.Range("D26")
$cell.Value = "'1.467"
$cell.Style = "Normal"
$ws.Range("E26").Value = '  +23.57%  '

$cell = $ws.Range("D27")
$cell.Value = "'2.435"
$cell.Style = "Normal"
$ws.Range("E27").Value = '  +0.53%  '

$cell = $ws.Range("D28")
$cell.Value = "'151.71"
$cell.Style = "Normal"
$ws.Range("E28").Value = '  +2.14%  '

$cell = $ws.Range("D29")
$cell.Value = "'19.56"
$cell.Style = "Normal"
$ws.Range("E29").Value = '  +4.00%  '

$ws.Range("D30").Value = '1.937.19'
$ws.Range("E30").Value = '  +3.82%  '

$cell = $ws.Range("D31")
$cell.Value = "'132.53"
$cell.Style = "Normal"
$ws.Range("E31").Value = '  +5.16%  '

$cell = $ws.Range("D32")
$cell.Value = "'4.110"
$cell.Style = "Normal"
$ws.Range("E32").Value = '  +0.85%  '

$cell = $ws.Range("D33")
$cell.Value = "'6.063"
$cell.Style = "Normal"
$ws.Range("E33").Value = '  +4.69%  '

$cell = $ws.Range("D34")
$cell.Value = "'0.08661"
$cell.Style = "Normal"
$ws.Range("E34").Value = '  +2.15%  '

$cell = $ws.Range("D35")
$cell.Value = "'1.695"
$cell.Style = "Normal"
$ws.Range("E35").Value = '  +2.60%  '

$cell = $ws.Range("D36")
$cell.Value = "'12.91"
$cell.Style = "Normal"
$ws.Range("E36").Value = '  +4.93%  '

$cell = $ws.Range("D37")
$cell.Value = "'5.441"
$cell.Style = "Normal"
$ws.Range("E37").Value = '  +4.79%  '

$cell = $ws.Range("D38")
$cell.Value = "'0.02359"
$cell.Style = "Normal"
$ws.Range("E38").Value = '  +4.04%  '

$cell = $ws.Range("D39")
$cell.Value = "'0.06277"
$cell.Style = "Normal"
$ws.Range("E39").Value = '  +3.96%  '

$cell = $ws.Range("D40")
$cell.Value = "'0.2177"
$cell.Style = "Normal"
$ws.Range("E40").Value = '  +3.64%  '

$cell = $ws.Range("D41")
$cell.Value = "'8.484"
$cell.Style = "Normal"
$ws.Range("E41").Value = '  +2.49%  '

$cell = $ws.Range("D42")
$cell.Value = "'1.223"
$cell.Style = "Normal"
$ws.Range("E42").Value = '  -0.39%  '

$cell = $ws.Range("D43")
$cell.Value = "'0.6274"
$cell.Style = "Normal"
$ws.Range("E43").Value = '  +4.92%  '

$cell = $ws.Range("D44")
$cell.Value = "'14.40"
$cell.Style = "Normal"
$ws.Range("E44").Value = '  +5.09%  '

$cell = $ws.Range("D45")
$cell.Value = "'0.9981"
$cell.Style = "Normal"

$cell = $ws.Range("D46")
$cell.Value = "'3.931"
$cell.Style = "Normal"
$ws.Range("E46").Value = '  +2.26%  '

$cell = $ws.Range("D47")
$cell.Value = "'0.6103"
$cell.Style = "Normal"
$ws.Range("E47").Value = '  +6.31%  '

$cell = $ws.Range("D48")
$cell.Value = "'129.34"
$cell.Style = "Normal"
$ws.Range("E48").Value = '  +2.85%  '

$cell = $ws.Range("D49")
$cell.Value = "'2.071"
$cell.Style = "Normal"
$ws.Range("E49").Value = '  +5.13%  '

$cell = $ws.Range("D50")
$cell.Value = "'0.07218"
$cell.Style = "Normal"
$ws.Range("E50").Value = '  +2.73%  '

$cell = $ws.Range("D51")
$cell.Value = "'77.85"
$cell.Style = "Normal"
$ws.Range("E51").Value = '  +4.03%  '
